$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on Price cells whose new values would otherwise be
# auto-detected as numbers by Excel, since the column stores text values
# (e.g. "594.66") formatted as plain text, not numeric cells.
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D51').NumberFormat = '@'

$ws.Range('D2').Value = '64.923.85'
$ws.Range('E2').Value = '  +0.42%  '
$ws.Range('D3').Value = '3.524.76'
$ws.Range('E3').Value = '  +3.10%  '
$ws.Range('E4').Value = '  +0.25%  '
$ws.Range('D5').Value = '594.66'
$ws.Range('E5').Value = '  +1.98%  '
$ws.Range('D6').Value = '136.45'
$ws.Range('E6').Value = '  +0.41%  '
$ws.Range('D7').Value = '3.527.82'
$ws.Range('E7').Value = '  +3.31%  '
$ws.Range('E8').Value = '  +0.15%  '
$ws.Range('E9').Value = '  +1.24%  '
$ws.Range('E10').Value = '  +2.00%  '
$ws.Range('E11').Value = '  -2.06%  '
$ws.Range('D12').Value = '0.381'
$ws.Range('E12').Value = '  +1.98%  '
$ws.Range('D13').Value = '4.124.17'
$ws.Range('E13').Value = '  +3.46%  '
$ws.Range('B14').Value = 'WrappedEther'
$ws.Range('C14').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D14').Value = '3.535.53'
$ws.Range('E14').Value = '  +2.74%  '
$ws.Range('B15').Value = 'Avalanche'
$ws.Range('C15').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D15').Value = '27.02'
$ws.Range('E15').Value = '  +2.71%  '
$ws.Range('B16').Value = 'ShibaInu'
$ws.Range('C16').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D16').Value = '0.0000180'
$ws.Range('E16').Value = '  +1.49%  '
$ws.Range('E17').Value = '  +1.44%  '
$ws.Range('D18').Value = '64.909.74'
$ws.Range('E18').Value = '  +0.54%  '
$ws.Range('D19').Value = '9.96'
$ws.Range('E19').Value = '  +4.33%  '
$ws.Range('D20').Value = '5.80'
$ws.Range('E20').Value = '  +0.13%  '
$ws.Range('D21').Value = '14.13'
$ws.Range('E21').Value = '  +4.69%  '
$ws.Range('D22').Value = '387.26'
$ws.Range('E22').Value = '  +1.16%  '
$ws.Range('D23').Value = '0.571'
$ws.Range('E23').Value = '  +4.19%  '
$ws.Range('D24').Value = '3.666.13'
$ws.Range('E24').Value = '  +3.46%  '
$ws.Range('D25').Value = '73.78'
$ws.Range('E25').Value = '  +2.69%  '
$ws.Range('D26').Value = '0.998'
$ws.Range('E26').Value = '  -0.18%  '
$ws.Range('D27').Value = '0.0000113'
$ws.Range('E27').Value = '  +7.59%  '
$ws.Range('E28').Value = '  +6.94%  '
$ws.Range('D29').Value = '1.00'
$ws.Range('E29').Value = '  +0.34%  '
$ws.Range('E30').Value = '  +3.38%  '
$ws.Range('E31').Value = '  +1.53%  '
$ws.Range('D32').Value = '3.538.81'
$ws.Range('E32').Value = '  +3.50%  '
$ws.Range('E33').Value = '  +0.02%  '
$ws.Range('D34').Value = '23.60'
$ws.Range('E34').Value = '  +2.97%  '
$ws.Range('D35').Value = '1.35'
$ws.Range('E35').Value = '  +14.09%  '
$ws.Range('E36').Value = '  +2.41%  '
$ws.Range('D37').Value = '170.22'
$ws.Range('E37').Value = '  +1.44%  '
$ws.Range('D38').Value = '1.54'
$ws.Range('E38').Value = '  +6.07%  '
$ws.Range('D39').Value = '6.79'
$ws.Range('E39').Value = '  +1.02%  '
$ws.Range('D40').Value = '4.92'
$ws.Range('E40').Value = '  +6.63%  '
$ws.Range('D41').Value = '0.0794'
$ws.Range('E41').Value = '  +5.49%  '
$ws.Range('D42').Value = '0.818'
$ws.Range('E42').Value = '  +1.41%  '
$ws.Range('E43').Value = '  +16.66%  '
$ws.Range('E44').Value = '  +0.41%  '
$ws.Range('D45').Value = '42.43'
$ws.Range('E45').Value = '  +1.24%  '
$ws.Range('D46').Value = '4.39'
$ws.Range('E46').Value = '  +2.46%  '
$ws.Range('E47').Value = '  +6.39%  '
$ws.Range('E48').Value = '  +3.43%  '
$ws.Range('D49').Value = '6.82'
$ws.Range('E49').Value = '  +6.03%  '
$ws.Range('D50').Value = '2.398.91'
$ws.Range('E50').Value = '  +11.07%  '
$ws.Range('D51').Value = '302.33'
$ws.Range('E51').Value = '  +10.50%  '
